# Add a new "model_V4" row (row 6) plus new columns I and J (Accuracy helper
# + Accuracy%) for every data row, matching the commit:
#   "updating our model to version 4 / adding more hidden layer ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6: model_V4 ---
$ws.Range("A6").Value = "model_V4"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1011
$ws.Range("D6").Value = 989
$ws.Range("E6").Value = 0

# Copy the (highlighted) format from E2 down onto E6 so it reuses the same
# cell style (s="1") instead of minting a new one.
$ws.Range("E2").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("F6").Formula = "=SUM(B6:E6)"

$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Formula = "=E6/F6"

# --- New helper columns I (E+D) and J (I/F) for every data row ---
$ws.Range("I2").Formula = "=E2+D2"
$ws.Range("I3").Formula = "=E3+D3"
$ws.Range("I4").Formula = "=E4+D4"
$ws.Range("I5").Formula = "=E5+D5"
$ws.Range("I6").Formula = "=E6+D6"

$ws.Range("J2").Formula = "=I2/F2"
$ws.Range("J2").NumberFormat = "0%"

# J3:J6 typed once and filled down together -> Excel records this as one
# shared formula group (J3 master, J4:J6 followers), matching the diff.
$ws.Range("J3:J6").Formula = "=I3/F3"
$ws.Range("J3:J6").NumberFormat = "0%"

# Update selection to match the diff (H7)
$ws.Range("H7").Select()
